$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the Check In Date / Check Out Date test values (row 2)
$ws.Range("G2").Value = "25/05/2016"
$ws.Range("H2").Value = "26/05/2016"

# Update the selected/active cell shown in the sheet view
$ws.Range("G12").Select()
